# Update "想去人数" (want-to-go count) figures on the sheets that list
# convention data. Both the "展览" sheet and the "全部类型" sheet carry the
# same rows, so the refreshed counts must be applied to both.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1773
    $ws.Range("F3").Value = 8093
    $ws.Range("F5").Value = 291
}
